$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 10 (Objetivos:) B/C text: replace misplaced professor name with the actual PT objective text ---
$ws.Range("B10").Value = "Apresentar os conceitos de spintrônica e as potenciais aplicações em computação quântica."
$ws.Range("C10").Value = "Apresentar os conceitos de spintrônica e as potenciais aplicações em computação quântica."

# --- Insert two new rows after row 12 (Docentes responsaveis:) for the professor list ---
$ws.Rows("13:14").Insert()
$ws.Range("A13:A14").Clear()

# Row 13: first professor (moved here from its old mis-placed spot)
$ws.Range("B13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C13").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("B3").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13").PasteSpecial(-4122)

# Row 14: second professor (new)
$ws.Range("B14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("C14").Value = "7290967 - Emerson Gonçalves de Melo"
$ws.Range("B3").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# --- Row 15 (was 13, Programa resumido:) B/C: replace stray date with the PT short-syllabus text ---
$ws.Range("B15").Value = "Introdução à nanotecnologia. Spintrônica de metais. Spintrônica de semicondutores. Dispositivos da spintrônica. Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos. Decoerência. Pontos quânticos. Transistor de Kane. Introdução a modelos da consciência: o cérebro é um computador quântico?"
$ws.Range("C15").Value = "Introdução à nanotecnologia. Spintrônica de metais. Spintrônica de semicondutores. Dispositivos da spintrônica. Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos. Decoerência. Pontos quânticos. Transistor de Kane. Introdução a modelos da consciência: o cérebro é um computador quântico?"

# --- Row 17 (was 15, Programa:) B/C: replace stray professor name with the PT full-syllabus text ---
$ws.Range("B17").Value = "Introdução à nanotecnologia.Spintrônica de metais. Spintrônica de semicondutores Dispositivos da spintrônica.Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos.Decoerência. Pontos quânticos. Transistor de Kane.Introdução a modelos da consciência: o cérebro é um computador quântico?"
$ws.Range("C17").Value = "Introdução à nanotecnologia.Spintrônica de metais. Spintrônica de semicondutores Dispositivos da spintrônica.Introdução à computação clássica. Introdução à computação quântica. Algoritmos quânticos.Decoerência. Pontos quânticos. Transistor de Kane.Introdução a modelos da consciência: o cérebro é um computador quântico?"

# --- Row 20 (was 18, Metodo:) B/C: replace stray professor name with the teaching-method text ---
$ws.Range("B20").Value = "Aulas expositivas, seminários e exercícios comentados."
$ws.Range("C20").Value = "Aulas expositivas, seminários e exercícios comentados."

# --- Row 21 (was 19, Criterio:) B/C: replace stray method text with the grading-criteria text ---
$ws.Range("B21").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."
$ws.Range("C21").Value = "Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2."

# --- Row 22 (was 20, Norma de recuperacao:) B/C: replace stray criteria text with the recovery-rule text ---
$ws.Range("B22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C22").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"

# --- Row 23 (was 21, Bibliografia:) B/C: replace stray recovery-rule text with the bibliography text ---
$ws.Range("B23").Value = "KITAEV, A. YU.; SHEN, A. H.; VYALVI, M. N. Classical and Quantum Computation, American Mathematical Society; 2002.`nBENENTI, G.; CASATI, G.; STRINI, G. Principles of Quantum Computation and Information, Vol. I: Basic Concepts, 2004.`nLO, H. K.; POPESCU, S.; SPILLER, T. Introduction to Quantum Computation and Information World Scientific Publishing Company, 2001."
$ws.Range("C23").Value = "KITAEV, A. YU.; SHEN, A. H.; VYALVI, M. N. Classical and Quantum Computation, American Mathematical Society; 2002.`nBENENTI, G.; CASATI, G.; STRINI, G. Principles of Quantum Computation and Information, Vol. I: Basic Concepts, 2004.`nLO, H. K.; POPESCU, S.; SPILLER, T. Introduction to Quantum Computation and Information World Scientific Publishing Company, 2001."

